# Update the crypto price/volume table with the latest scraped values.
# GitHub Actions scheduled refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "26.80", "1.97").
# Force the whole data range to Text format first so Excel doesn't coerce
# these into real numbers (which would also drop trailing zeros), then
# restore the default "Normal" style afterwards so no stray formatting is
# left behind on cells we touch.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# --- Row 2: Bitcoin ---
$ws.Range("D2").Value = "29.589.40"
$ws.Range("E2").Value = "  +1.42%  "

# --- Row 3: Ethereum ---
$ws.Range("D3").Value = "1.600.51"
$ws.Range("E3").Value = "  +1.37%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  +0.37%  "

# --- Row 5: BNB ---
$ws.Range("D5").Value = "212.34"
$ws.Range("E5").Value = "  +0.15%  "

# --- Row 6: XRP ---
$ws.Range("E6").Value = "  +0.55%  "

# --- Row 7: USDC ---
$ws.Range("E7").Value = "  +0.24%  "

# --- Row 8: Solana ---
$ws.Range("D8").Value = "26.80"
$ws.Range("E8").Value = "  +1.56%  "

# --- Row 9: Cardano ---
$ws.Range("E9").Value = "  +1.12%  "

# --- Row 10: Dogecoin ---
$ws.Range("E10").Value = "  +1.27%  "

# --- Row 11: TRON ---
$ws.Range("E11").Value = "  +1.05%  "

# --- Row 12: WrappedliquidstakedEther2.0 ---
$ws.Range("D12").Value = "1.830.27"
$ws.Range("E12").Value = "  +1.45%  "

# --- Row 13: WrappedEther ---
$ws.Range("D13").Value = "1.596.55"
$ws.Range("E13").Value = "  +1.18%  "

# --- Row 14/15: WrappedBTC and Polygon swapped rank order ---
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "0.539"
$ws.Range("E14").Value = "  +2.98%  "

$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "29.598.74"
$ws.Range("E15").Value = "  +1.34%  "

# --- Row 16: Polkadot ---
$ws.Range("E16").Value = "  +0.71%  "

# --- Row 17: Litecoin ---
$ws.Range("D17").Value = "63.82"
$ws.Range("E17").Value = "  +2.40%  "

# --- Row 18: BitcoinCash ---
$ws.Range("D18").Value = "241.45"
$ws.Range("E18").Value = "  +1.77%  "

# --- Row 19 ---
$ws.Range("E19").Value = "  +2.37%  "

# --- Row 20 ---
$ws.Range("D20").Value = "0.0₃0694"
$ws.Range("E20").Value = "  +0.40%  "

# --- Row 21 ---
$ws.Range("E21").Value = "  +0.60%  "

# --- Row 22 ---
$ws.Range("E22").Value = "  -0.33%  "

# --- Row 23 ---
$ws.Range("E23").Value = "  +0.44%  "

# --- Row 24 ---
$ws.Range("E24").Value = "  -1.20%  "

# --- Row 25 ---
$ws.Range("D25").Value = "155.10"
$ws.Range("E25").Value = "  +0.99%  "

# --- Row 26 ---
$ws.Range("E26").Value = "  +1.23%  "

# --- Row 27 ---
$ws.Range("E27").Value = "  +0.77%  "

# --- Row 28 ---
$ws.Range("E28").Value = "  +1.06%  "

# --- Row 29 ---
$ws.Range("E29").Value = "  +0.40%  "

# --- Row 30 ---
$ws.Range("E30").Value = "  +2.51%  "

# --- Row 31 ---
$ws.Range("E31").Value = "  +0.30%  "

# --- Row 32 ---
$ws.Range("D32").Value = "3.22"
$ws.Range("E32").Value = "  +0.27%  "

# --- Row 33 ---
$ws.Range("D33").Value = "3.16"
$ws.Range("E33").Value = "  +2.86%  "

# --- Row 34 ---
$ws.Range("D34").Value = "1.423.44"
$ws.Range("E34").Value = "  +0.00%  "

# --- Row 35 ---
$ws.Range("D35").Value = "1.55"
$ws.Range("E35").Value = "  +2.37%  "

# --- Row 36 ---
$ws.Range("E36").Value = "  +4.35%  "

# --- Row 37 ---
$ws.Range("D37").Value = "1.03"
$ws.Range("E37").Value = "  -1.74%  "

# --- Row 38 ---
$ws.Range("D38").Value = "2.29"
$ws.Range("E38").Value = "  +0.32%  "

# --- Row 39 ---
$ws.Range("D39").Value = "0.0168"
$ws.Range("E39").Value = "  +3.09%  "

# --- Row 40 ---
$ws.Range("D40").Value = "0.543"
$ws.Range("E40").Value = "  +2.95%  "

# --- Row 41/42/43: RenderToken, Kaspa, BitcoinSV rank order shuffled ---
$ws.Range("B41").Value = "BitcoinSV"
$ws.Range("C41").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D41").Value = "55.32"
$ws.Range("E41").Value = "  +4.81%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "1.97"
$ws.Range("E42").Value = "  +0.13%  "

$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "0.0495"
$ws.Range("E43").Value = "  +5.84%  "

# --- Row 44: ARBITRUM ---
$ws.Range("D44").Value = "0.809"
$ws.Range("E44").Value = "  +2.61%  "

# --- Row 45: PaxDollar ---
$ws.Range("E45").Value = "  +0.41%  "

# --- Row 46: WEMIXToken ---
$ws.Range("D46").Value = "0.989"
$ws.Range("E46").Value = "  +16.77%  "

# --- Row 47: Aave ---
$ws.Range("D47").Value = "66.32"
$ws.Range("E47").Value = "  +2.68%  "

# --- Row 48: FraxShare ---
$ws.Range("E48").Value = "  -0.58%  "

# --- Row 49: RocketPoolETH ---
$ws.Range("D49").Value = "1.740.53"
$ws.Range("E49").Value = "  +1.39%  "

# --- Row 50: Quant ---
$ws.Range("D50").Value = "85.97"
$ws.Range("E50").Value = "  +0.24%  "

# --- Row 51: BabyDogeCoin ---
$ws.Range("D51").Value = "0.0₆0105"
$ws.Range("E51").Value = "  +2.06%  "

# Drop the temporary text-number format we applied, restoring the original
# default cell style so only the cell values themselves changed.
$priceRange.Style = "Normal"
